$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H5").Value = 228.7
$ws.Range("I5").Value = 250.57143
$ws.Range("K5").Value = 250.57143
$ws.Range("M5").Value = -135.57143
$ws.Range("H106").Value = 48276.5
$ws.Range("I106").Value = 51724.582
$ws.Range("J106").Value = 6899.5
$ws.Range("K106").Value = 51724.582
$ws.Range("L106").Value = 6899.5
$ws.Range("M106").Value = -51093.582
$ws.Range("N106").Value = -8161.5
$ws.Range("H132").Value = 5402.982
$ws.Range("I132").Value = 4212.8887
$ws.Range("K132").Value = 12638.6661
$ws.Range("M132").Value = -10108.6661
$ws.Range("H137").Value = 2481.487
$ws.Range("I137").Value = 1901.56
$ws.Range("J137").Value = 3517.0715
$ws.Range("K137").Value = 5704.68
$ws.Range("L137").Value = 10551.2145
$ws.Range("M137").Value = -3154.68
$ws.Range("N137").Value = -15651.2145
$ws.Range("H138").Value = 194406.42
$ws.Range("J138").Value = 530096.9399999999
$ws.Range("L138").Value = 1590290.82
$ws.Range("N138").Value = -1600570.82

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H5").Value = 3747.4
$ws.Range("I5").Value = 3997.2144
$ws.Range("K5").Value = 3997.2144
$ws.Range("M5").Value = -3885.2144
$ws.Range("H32").Value = 11570.023
$ws.Range("I32").Value = 11570.023
$ws.Range("J32").Value = 0
$ws.Range("K32").Value = 11570.023
$ws.Range("L32").Value = 0
$ws.Range("M32").Value = -11283.023
$ws.Range("N32").ClearContents()
$ws.Range("H34").Value = 30025
$ws.Range("I34").Value = 30025
$ws.Range("K34").Value = 30025
$ws.Range("M34").Value = -29754
$ws.Range("H40").Value = 30028
$ws.Range("I40").Value = 30028
$ws.Range("K40").Value = 30028
$ws.Range("M40").Value = -29852
$ws.Range("H74").Value = 1427.0526
$ws.Range("I74").Value = 1365.7059
$ws.Range("J74").Value = 1948.5
$ws.Range("K74").Value = 1365.7059
$ws.Range("L74").Value = 1948.5
$ws.Range("M74").Value = -491.7058999999999
$ws.Range("N74").Value = -3696.5
$ws.Range("H77").Value = 1427.0526
$ws.Range("I77").Value = 1365.7059
$ws.Range("J77").Value = 1948.5
$ws.Range("K77").Value = 6828.5295
$ws.Range("L77").Value = 9742.5
$ws.Range("M77").Value = -2460.5295
$ws.Range("N77").Value = -18478.5
$ws.Range("H92").Value = 90000
$ws.Range("J92").Value = 90000
$ws.Range("L92").Value = 90000
$ws.Range("N92").Value = -94992
$ws.Range("H97").Value = 1020.6
$ws.Range("J97").Value = 1377
$ws.Range("L97").Value = 1377
$ws.Range("N97").Value = -2369
$ws.Range("H132").Value = 2783.8909
$ws.Range("I132").Value = 1279.2709
$ws.Range("K132").Value = 3837.8127
$ws.Range("M132").Value = -1307.8127

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H4").Value = 3747.4
$ws.Range("I4").Value = 3997.2144
$ws.Range("K4").Value = 3997.2144
$ws.Range("M4").Value = -3882.2144

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 431.8846
$ws.Range("I22").Value = 441.2
$ws.Range("K22").Value = 441.2
$ws.Range("M22").Value = -91.19999999999999
$ws.Range("H47").Value = 24250
$ws.Range("J47").Value = 0
$ws.Range("L47").Value = 0
$ws.Range("N47").ClearContents()
$ws.Range("H58").Value = 3553.121
$ws.Range("I58").Value = 3573.28
$ws.Range("J58").Value = 3490.125
$ws.Range("K58").Value = 3573.28
$ws.Range("L58").Value = 3490.125
$ws.Range("M58").Value = -3370.28
$ws.Range("N58").Value = -3896.125
$ws.Range("H99").Value = 5455.5713
$ws.Range("I99").Value = 5287.75
$ws.Range("K99").Value = 5287.75
$ws.Range("M99").Value = -3789.75
$ws.Range("H105").Value = 3000
$ws.Range("I105").Value = 2338
$ws.Range("J105").Value = 4986
$ws.Range("K105").Value = 2338
$ws.Range("L105").Value = 4986
$ws.Range("M105").Value = -591
$ws.Range("N105").Value = -8480
$ws.Range("H126").Value = 5455.5713
$ws.Range("I126").Value = 5287.75
$ws.Range("K126").Value = 15863.25
$ws.Range("M126").Value = -13393.25
$ws.Range("H134").Value = 1840.4
$ws.Range("I134").Value = 1918.7646
$ws.Range("K134").Value = 5756.293799999999
$ws.Range("M134").Value = -3221.293799999999
$ws.Range("H136").Value = 3553.121
$ws.Range("I136").Value = 3573.28
$ws.Range("J136").Value = 3490.125
$ws.Range("K136").Value = 10719.84
$ws.Range("L136").Value = 10470.375
$ws.Range("M136").Value = -8169.84
$ws.Range("N136").Value = -15570.375

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H12").Value = 505.7037
$ws.Range("J12").Value = 363.70587
$ws.Range("L12").Value = 1091.11761
$ws.Range("N12").Value = -1437.11761
$ws.Range("H34").Value = 72.5
$ws.Range("I34").Value = 72.5
$ws.Range("K34").Value = 217.5
$ws.Range("M34").Value = -133.5
$ws.Range("H39").Value = 1958.2
$ws.Range("J39").Value = 3450
$ws.Range("L39").Value = 10350
$ws.Range("N39").Value = -10938
$ws.Range("H124").Value = 15999
$ws.Range("J124").Value = 15999
$ws.Range("L124").Value = 47997
$ws.Range("N124").Value = -57817
$ws.Range("H129").Value = 2527.889
$ws.Range("I129").Value = 1120.5
$ws.Range("K129").Value = 3361.5
$ws.Range("M129").Value = 1638.5
$ws.Range("H131").Value = 1531.8379
$ws.Range("J131").Value = 1784.6786
$ws.Range("L131").Value = 5354.0358
$ws.Range("N131").Value = -15434.0358

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 282
$ws.Range("J2").Value = 298
$ws.Range("L2").Value = 298
$ws.Range("N2").Value = -524
$ws.Range("H18").Value = 0
$ws.Range("J18").Value = 0
$ws.Range("L18").Value = 0
$ws.Range("N18").ClearContents()
$ws.Range("H33").Value = 29999.5
$ws.Range("J33").Value = 29999.5
$ws.Range("L33").Value = 29999.5
$ws.Range("N33").Value = -30503.5
$ws.Range("H107").Value = 821.1
$ws.Range("J107").Value = 853.4666999999999
$ws.Range("L107").Value = 853.4666999999999
$ws.Range("N107").Value = -4693.4667
$ws.Range("H123").Value = 55000
$ws.Range("J123").Value = 55000
$ws.Range("L123").Value = 55000
$ws.Range("N123").Value = -59900

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 2673.1875
$ws.Range("J22").Value = 3398.4546
$ws.Range("L22").Value = 3398.4546
$ws.Range("N22").Value = -3988.4546
$ws.Range("H27").Value = 2673.1875
$ws.Range("J27").Value = 3398.4546
$ws.Range("L27").Value = 3398.4546
$ws.Range("N27").Value = -3612.4546
$ws.Range("H41").Value = 30033
$ws.Range("I41").Value = 30033
$ws.Range("K41").Value = 30033
$ws.Range("M41").Value = -29595
$ws.Range("H46").Value = 21035.133
$ws.Range("I46").Value = 5694.25
$ws.Range("J46").Value = 26613.637
$ws.Range("K46").Value = 5694.25
$ws.Range("L46").Value = 26613.637
$ws.Range("M46").Value = -5506.25
$ws.Range("N46").Value = -26989.637
$ws.Range("H47").Value = 0
$ws.Range("I47").Value = 0
$ws.Range("J47").Value = 0
$ws.Range("K47").Value = 0
$ws.Range("L47").Value = 0
$ws.Range("M47").ClearContents()
$ws.Range("N47").ClearContents()
$ws.Range("H52").Value = 0
$ws.Range("I52").Value = 0
$ws.Range("J52").Value = 0
$ws.Range("K52").Value = 0
$ws.Range("L52").Value = 0
$ws.Range("M52").ClearContents()
$ws.Range("N52").ClearContents()
$ws.Range("H100").Value = 4065.1904
$ws.Range("I100").Value = 4060.4614
$ws.Range("J100").Value = 4072.875
$ws.Range("K100").Value = 4060.4614
$ws.Range("L100").Value = 4072.875
$ws.Range("M100").Value = -3519.4614
$ws.Range("N100").Value = -5154.875
$ws.Range("H122").Value = 5169.5625
$ws.Range("I122").Value = 3427.1765
$ws.Range("K122").Value = 10281.5295
$ws.Range("M122").Value = -7831.529500000001
$ws.Range("H133").Value = 89998
$ws.Range("J133").Value = 89998
$ws.Range("L133").Value = 89998
$ws.Range("N133").Value = -95058

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H49").Value = 31841.715
$ws.Range("I49").Value = 32590.8
$ws.Range("J49").Value = 29969
$ws.Range("K49").Value = 32590.8
$ws.Range("L49").Value = 29969
$ws.Range("M49").Value = -32360.8
$ws.Range("N49").Value = -30429
$ws.Range("H54").Value = 11820.444
$ws.Range("I54").Value = 9197.714
$ws.Range("J54").Value = 21000
$ws.Range("K54").Value = 9197.714
$ws.Range("L54").Value = 21000
$ws.Range("M54").Value = -8677.714
$ws.Range("N54").Value = -22040
$ws.Range("H62").Value = 75867.91
$ws.Range("I62").Value = 82654.7
$ws.Range("K62").Value = 82654.7
$ws.Range("M62").Value = -82030.7
$ws.Range("H65").Value = 75867.91
$ws.Range("I65").Value = 82654.7
$ws.Range("K65").Value = 413273.5
$ws.Range("M65").Value = -410153.5
$ws.Range("H122").Value = 2884.276
$ws.Range("I122").Value = 2869
$ws.Range("K122").Value = 8607
$ws.Range("M122").Value = -6157
$ws.Range("H133").Value = 52749.25
$ws.Range("J133").Value = 52749.25
$ws.Range("L133").Value = 52749.25
$ws.Range("N133").Value = -62869.25
